# Update gh-pages to output generated at 456a3b4
# Applies to sheets "展览" (index 1) and "全部类型" (index 4):
#  - bump the "想去人数" (F column) counts for the existing 5 events
#  - insert a new event "南宁·熊喵M动漫嘉年华·万圣派对" before the last
#    existing event, shifting the last event down one row and bumping
#    its "想去人数" count as well

$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Bump want-to-go counts on rows 2-6
    $ws.Range("F2").Value = 45
    $ws.Range("F3").Value = 121
    $ws.Range("F4").Value = 156
    $ws.Range("F5").Value = 3129
    $ws.Range("F6").Value = 315

    # Last existing data row (index of the 万圣漫控嘉年华10 row) before insertion
    $lastRow = $ws.UsedRange.Rows.Count

    # Remember the existing last row's values so we can move them down
    $oldA = $ws.Cells.Item($lastRow, 1).Value()
    $oldB = $ws.Cells.Item($lastRow, 2).Value()
    $oldC = $ws.Cells.Item($lastRow, 3).Value()
    $oldD = $ws.Cells.Item($lastRow, 4).Value()
    $oldE = $ws.Cells.Item($lastRow, 5).Value()
    $oldF = $ws.Cells.Item($lastRow, 6).Value()
    $oldG = $ws.Cells.Item($lastRow, 7).Value()
    $oldH = $ws.Cells.Item($lastRow, 8).Value()
    $oldI = $ws.Cells.Item($lastRow, 9).Value()

    # Insert a new blank row right after the old last row, it will hold
    # the old last row's data (shifted down) while the old row position
    # gets the new event's data
    $newRow = $lastRow + 1
    $ws.Rows.Item($newRow).Insert()

    # Copy the formatting (bordered/bold/centered style) of the old last
    # row's sequence-number cell onto the new row's sequence-number cell
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # Write moved-down data for the previously-last event into the new row,
    # bumping its want-to-go count by 1 and its sequence number by 1.
    # The date column (oldB, e.g. "2024-11-02") looks like a real date to
    # Excel's auto-detection, so force it to be stored as plain text.
    $ws.Cells.Item($newRow, 1).Value = $oldA + 1
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = $oldB
    $ws.Cells.Item($newRow, 2).Style = "Normal"
    $ws.Cells.Item($newRow, 3).Value = $oldC
    $ws.Cells.Item($newRow, 4).Value = $oldD
    $ws.Cells.Item($newRow, 5).Value = $oldE
    $ws.Cells.Item($newRow, 6).Value = $oldF + 1
    $ws.Cells.Item($newRow, 7).Value = $oldG
    $ws.Cells.Item($newRow, 8).Value = $oldH
    $ws.Cells.Item($newRow, 9).Value = $oldI

    # Overwrite the old last row's cells with the new event's data
    # (its sequence number in column A stays the same).
    # The date column looks like a real date to Excel's auto-detection,
    # so force it to be stored as plain text instead of a date serial.
    $ws.Cells.Item($lastRow, 2).NumberFormat = "@"
    $ws.Cells.Item($lastRow, 2).Value = "2024-10-26"
    $ws.Cells.Item($lastRow, 2).Style = "Normal"
    $ws.Cells.Item($lastRow, 3).Value = "南宁·熊喵M动漫嘉年华·万圣派对"
    $ws.Cells.Item($lastRow, 4).Value = "亭洪路45号 百益上河城"
    $ws.Cells.Item($lastRow, 5).Value = "2024.10.26 11:00-10.27 21:00"
    $ws.Cells.Item($lastRow, 6).Value = 2
    $ws.Cells.Item($lastRow, 7).Value = 60
    $ws.Cells.Item($lastRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91894"
    $ws.Cells.Item($lastRow, 9).Value = "//i2.hdslb.com/bfs/openplatform/202409/NzSGhcoK1725123076559.jpeg"
}
